# D3_source_population.xlsx — rename two header columns and move the
# active selection, matching the author's "updated figures and tables" edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row renames: "birthdate" -> "birth_date", "sex" -> "gender"
$ws.Range("B1").Value = "birth_date"
$ws.Range("C1").Value = "gender"

# Move the active selection from E5 to C1
$ws.Range("C1").Select()
